$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: advance the date in A1 by one day
$ws.Range("A1").Value = 45309

# Step 2: update the price list (column D) values
$ws.Range("D26").Value = 1246
$ws.Range("D27").Value = 1246
$ws.Range("D28").Value = 1467
$ws.Range("D29").Value = 2536
$ws.Range("D30").Value = 3128
$ws.Range("D31").Value = 4150
$ws.Range("D32").Value = 5903
$ws.Range("D33").Value = 9740

$ws.Range("D35").Value = 1840
$ws.Range("D36").Value = 2082
$ws.Range("D37").Value = 2657
$ws.Range("D38").Value = 4414
$ws.Range("D39").Value = 5725
$ws.Range("D40").Value = 7299
$ws.Range("D41").Value = 10263
$ws.Range("D42").Value = 15834
